# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibitions) sheet and the corresponding rows on the
# "全部类型" (All types) combined sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 3055
$wsExpo.Range("F4").Value  = 107
$wsExpo.Range("F5").Value  = 6796
$wsExpo.Range("F6").Value  = 1810
$wsExpo.Range("F7").Value  = 41
$wsExpo.Range("F9").Value  = 32
$wsExpo.Range("F10").Value = 64
$wsExpo.Range("F11").Value = 133
$wsExpo.Range("F12").Value = 150
$wsExpo.Range("F13").Value = 28

# --- Sheet "全部类型" (sheet4), offset by +1 row vs. "展览" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 3055
$wsAll.Range("F5").Value  = 107
$wsAll.Range("F6").Value  = 6796
$wsAll.Range("F7").Value  = 1810
$wsAll.Range("F8").Value  = 41
$wsAll.Range("F10").Value = 32
$wsAll.Range("F11").Value = 64
$wsAll.Range("F12").Value = 133
$wsAll.Range("F13").Value = 150
$wsAll.Range("F14").Value = 28
